$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header H1 ("Save") - reuse the same formatting as the other
# header cells (bold, centered, bordered) by copying format from G1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data column values (plain numbers, default style)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0

$excel.CutCopyMode = 0
